# "Fix Melbourne sparse data"
#
# Updates the sparse-data table on the "Melbourne" sheet with corrected
# numeric values (several previously-blank F/G/H "_dev" cells are now
# populated, and a number of existing p10/p50/p90 mean+dev values change),
# sets the Melbourne sheet's page setup (paper size / orientation), clears
# the stray custom row height on the "Stuttgart" sheet's row 4 (back to the
# sheet default), and finally leaves "Melbourne" as the active sheet with
# A15 selected (matching the saved workbook view state).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Stuttgart: row 4 had an explicit 23.4pt row height; restore the sheet's
# default row height (drops the custom height instead of merely matching
# the default value numerically).
# ---------------------------------------------------------------------
$stuttgart = $wb.Worksheets.Item("Stuttgart")
$stuttgart.Rows.Item(4).AutoFit()

# ---------------------------------------------------------------------
# Melbourne: corrected sparse-data values.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Melbourne")
$ws.Activate()

# Row 2 (p10_mean / p50_mean / p90_mean + newly-populated *_dev columns)
$ws.Range("B2").Value = 111709.5425
$ws.Range("C2").Value = 112375.7389
$ws.Range("D2").Value = 112880.9439
$ws.Range("F2").Value = 354.74570590000002
$ws.Range("G2").Value = 275.12807959999998
$ws.Range("H2").Value = 306.4205336

# Row 3
$ws.Range("B3").Value = 104734.02650000001
$ws.Range("C3").Value = 105222.6039
$ws.Range("D3").Value = 105793.79580000001
$ws.Range("F3").Value = 236.0989309
$ws.Range("G3").Value = 187.6438417
$ws.Range("H3").Value = 278.69779690000001

# Row 5
$ws.Range("B5").Value = 15984.173059999999
$ws.Range("C5").Value = 17092.620699999999
$ws.Range("D5").Value = 17897.945329999999
$ws.Range("F5").Value = 74.064079370000002
$ws.Range("G5").Value = 322.99627889999999
$ws.Range("H5").Value = 119.69559340000001

# Row 7
$ws.Range("B7").Value = 0.00040232799999999997
$ws.Range("C7").Value = 0.00058297100000000005
$ws.Range("D7").Value = 0.000750302
$ws.Range("F7").Value = 0.000018300000000000001
$ws.Range("G7").Value = 0.000038399999999999998
$ws.Range("H7").Value = 0.000060399999999999998

# Row 8 (B/C/D already 0 - only F/G/H newly populated)
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 0

# Row 9
$ws.Range("B9").Value = 0.0036448349999999999
$ws.Range("C9").Value = 0.0038843380000000002
$ws.Range("D9").Value = 0.0042781520000000003

# Row 10 (C10 unchanged; F/G/H newly populated with the old B/C/D values)
$ws.Range("B10").Value = 0.000666583
$ws.Range("D10").Value = 0.000857116
$ws.Range("F10").Value = 0.00066722099999999996
$ws.Range("G10").Value = 0.00075796200000000005
$ws.Range("H10").Value = 0.00085354300000000001

# Row 12
$ws.Range("B12").Value = 0.000000000000000031599999999999998
$ws.Range("C12").Value = 0.0000000000000000339
$ws.Range("D12").Value = 0.000000000000000038099999999999999

# Row 14
$ws.Range("B14").Value = 0.0015702730000000001
$ws.Range("C14").Value = 0.001611407
$ws.Range("D14").Value = 0.001619606

# Row 15
$ws.Range("B15").Value = 0.0000025799999999999999
$ws.Range("C15").Value = 0.0000038099999999999999
$ws.Range("D15").Value = 0.0000234

# Row 17
$ws.Range("F17").Value = 125.9389721
$ws.Range("G17").Value = 152.95976490000001
$ws.Range("H17").Value = 268.41984600000001

# ---------------------------------------------------------------------
# Page setup for the Melbourne sheet (set via Print Setup dialog).
# ---------------------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# ---------------------------------------------------------------------
# Final view state: Melbourne tab active, A15 selected.
# ---------------------------------------------------------------------
$ws.Range("A15").Select() | Out-Null
